# Quarterly database update + column (quarter) rollover:
# drop the oldest quarter label/column of data and append the new quarter
# ("فصل سوم منتهی به 1401/12") with freshly computed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Quarter header labels (row 8 and row 24), columns E:N ---
$quarters = @(
    "فصل دوم منتهی به 1399/09",
    "فصل سوم منتهی به 1399/12",
    "فصل چهارم منتهی به 1400/03",
    "فصل اول منتهی به 1400/06",
    "فصل دوم منتهی به 1400/09",
    "فصل سوم منتهی به 1400/12",
    "فصل چهارم منتهی به 1401/03",
    "فصل اول منتهی به 1401/06",
    "فصل دوم منتهی به 1401/09",
    "فصل سوم منتهی به 1401/12"
)

$cols = @("E","F","G","H","I","J","K","L","M","N")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $quarters[$i]
    $ws.Range($cols[$i] + "24").Value = $quarters[$i]
}

# --- Data rows: shift one quarter left, append newly computed quarter value ---

# Row 12 - هزینه حمل و نقل و انتقال
$row12 = @(1506, 489, 1095, 908, 839, 1135, 4931, 5596, 17166, -10962)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "12").Value = $row12[$i] }

# Row 16 - هزینه انرژی (آب، برق، گاز و سوخت)
$row16 = @(565, 806, 946, 1281, 1383, 1457, 1653, 5627, 5983, 3322)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "16").Value = $row16[$i] }

# Row 17 - هزینه استهلاک
$row17 = @(14800, 17700, 22742, 19388, 26332, 19062, 43747, 42848, 24081, 40624)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "17").Value = $row17[$i] }

# Row 19 - هزینه مطالبات مشکوک الوصول
$row19 = @(27301, 15113, 6345, 13361, 39843, 31475, 15016, 31189, 73837, 91786)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "19").Value = $row19[$i] }

# Row 20 - جمع
$row20 = @(44172, 34108, 31128, 34938, 68397, 53129, 65347, 85260, 121067, 124770)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "20").Value = $row20[$i] }

# Row 26 - تعداد پرسنل غیر تولیدی شرکت
$row26 = @(41, 38, 39, 43, 43, 45, 48, 48, 48, 46)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "26").Value = $row26[$i] }

# Row 27 - تعداد پرسنل تولیدی شرکت
$row27 = @(566, 567, 581, 583, 607, 637, 670, 689, 714, 705)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "27").Value = $row27[$i] }
